# Generate Report for Handoff
# - Update Priority from "low" to "ht" for the 4 non-md-source rows (rows 4-7)
#   on both the "zh-cn" and "de-de" worksheets.
# - Refresh the "Latest Handoff Datetime" for the same zh-cn rows.
# - Refresh the "Latest HO Xliff Generate Date" column on the Overview sheet
#   for the corresponding 4 rows.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-08-31 16:35:35"
$wsZhCn.Range("H5").Value = "2016-08-31 16:35:35"
$wsZhCn.Range("H6").Value = "2016-08-31 16:35:35"
$wsZhCn.Range("H7").Value = "2016-08-31 16:35:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 16:35:39"
$wsOverview.Range("G5").Value = "2016-08-31 16:35:39"
$wsOverview.Range("G6").Value = "2016-08-31 16:35:39"
$wsOverview.Range("G7").Value = "2016-08-31 16:35:39"

# The "Latest Handoff Datetime" cells on de-de share the very same backing
# string as Overview's "Latest HO Xliff Generate Date" column, so refresh
# them to the same new timestamp to keep both in sync.
$wsDeDe.Range("H4").Value = "2016-08-31 16:35:39"
$wsDeDe.Range("H5").Value = "2016-08-31 16:35:39"
$wsDeDe.Range("H6").Value = "2016-08-31 16:35:39"
$wsDeDe.Range("H7").Value = "2016-08-31 16:35:39"
